$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param([string]$CellRef, [string]$TextValue)
    $c = $ws.Range($CellRef)
    # Force text type so numeric-looking strings (e.g. "1.00", "58.041.83")
    # are stored as text, matching the source data (inline/shared strings),
    # not auto-coerced into numbers by Excel.
    $c.NumberFormat = "@"
    $c.Value = $TextValue
    # Drop the temporary text-format style so the cell ends up with the
    # same (default) style it had before, just like the other text cells.
    $c.ClearFormats()
}

Set-TextValue 'D2' '58.041.83'
Set-TextValue 'E2' '  +2.46%  '
Set-TextValue 'D3' '3.065.73'
Set-TextValue 'E3' '  +2.65%  '
Set-TextValue 'E4' '  -0.09%  '
Set-TextValue 'D5' '528.28'
Set-TextValue 'E5' '  +6.21%  '
Set-TextValue 'D6' '143.89'
Set-TextValue 'E6' '  +6.55%  '
Set-TextValue 'E7' '  -0.04%  '
Set-TextValue 'D8' '0.449'
Set-TextValue 'E8' '  +5.34%  '
Set-TextValue 'D9' '7.65'
Set-TextValue 'E9' '  +5.39%  '
Set-TextValue 'E10' '  +7.14%  '
Set-TextValue 'E11' '  +6.03%  '
Set-TextValue 'E12' '  +2.05%  '
Set-TextValue 'D13' '3.588.98'
Set-TextValue 'E13' '  +2.58%  '
Set-TextValue 'D14' '27.40'
Set-TextValue 'E14' '  +7.99%  '
Set-TextValue 'D15' '0.0000173'
Set-TextValue 'E15' '  +16.50%  '
Set-TextValue 'D16' '57.988.04'
Set-TextValue 'E16' '  +2.37%  '
Set-TextValue 'E17' '  +8.54%  '
Set-TextValue 'D18' '3.059.42'
Set-TextValue 'E18' '  +2.59%  '
Set-TextValue 'D19' '13.24'
Set-TextValue 'E19' '  +7.29%  '
Set-TextValue 'D20' '8.22'
Set-TextValue 'E20' '  +5.75%  '
Set-TextValue 'D21' '343.02'
Set-TextValue 'E21' '  +5.03%  '
Set-TextValue 'E22' '  +0.08%  '
Set-TextValue 'D23' '5.69'
Set-TextValue 'E23' '  -0.85%  '
Set-TextValue 'E24' '  +7.84%  '
Set-TextValue 'D25' '65.48'
Set-TextValue 'E25' '  +5.99%  '
Set-TextValue 'E26' '  +9.51%  '
Set-TextValue 'E27' '  +4.99%  '
Set-TextValue 'D28' '0.999'
Set-TextValue 'E28' '  -0.20%  '
Set-TextValue 'E29' '  +9.16%  '
Set-TextValue 'E30' '  +9.87%  '
Set-TextValue 'E31' '  +7.22%  '
Set-TextValue 'E32' '  +5.68%  '
Set-TextValue 'D33' '21.29'
Set-TextValue 'E33' '  +4.05%  '
Set-TextValue 'E34' '  +8.56%  '
Set-TextValue 'D35' '157.87'
Set-TextValue 'E35' '  +3.19%  '
Set-TextValue 'E36' '  +7.20%  '
Set-TextValue 'E37' '  +4.38%  '
Set-TextValue 'D38' '26.22'
Set-TextValue 'E38' '  +12.84%  '
Set-TextValue 'E39' '  +4.63%  '
Set-TextValue 'D40' '3.098.63'
Set-TextValue 'D41' '37.83'
Set-TextValue 'E41' '  +3.73%  '
Set-TextValue 'D42' '3.96'
Set-TextValue 'E42' '  +11.42%  '
Set-TextValue 'E43' '  +5.04%  '
Set-TextValue 'B44' 'Maker'
Set-TextValue 'C44' 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue 'D44' '2.348.77'
Set-TextValue 'E44' '  +5.80%  '
Set-TextValue 'D45' '0.668'
Set-TextValue 'E45' '  +4.41%  '
Set-TextValue 'B46' 'Stacks'
Set-TextValue 'C46' 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue 'D46' '1.48'
Set-TextValue 'E46' '  +6.02%  '
Set-TextValue 'B47' 'FirstDigitalUSD'
Set-TextValue 'C47' 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue 'D47' '1.00'
Set-TextValue 'E47' '  -0.08%  '
Set-TextValue 'D48' '6.14'
Set-TextValue 'E48' '  +7.15%  '
Set-TextValue 'E49' '  +4.04%  '
Set-TextValue 'D50' '0.0246'
Set-TextValue 'E50' '  +4.43%  '
Set-TextValue 'D51' '20.30'
Set-TextValue 'E51' '  +6.51%  '
